# Auto-generated Excel COM-interop script
# Applies the scheduled-runner profit recalculation updates to the Leve profit
# columns (H:N) across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 33: Glazed and Confused | Clear Glass Lens
$ws.Range("H33").Value = 126.44444
$ws.Range("I33").Value = 134
$ws.Range("K33").Value = 134
$ws.Range("M33").Value = 95

# Row 38: Just Give Him a Serum | Hi-Potion of Strength
$ws.Range("H38").Value = 1820.2

# Row 55: A Real Smooth Move | Lanolin
$ws.Range("H55").Value = 780.6
$ws.Range("I55").Value = 777.1
$ws.Range("J55").Value = 787.6
$ws.Range("K55").Value = 777.1
$ws.Range("L55").Value = 787.6
$ws.Range("M55").Value = -563.1
$ws.Range("N55").Value = -1215.6

# Row 69: Steeling the Knife, Steeling the Mind | Grade 1 Mind Dissolvent
$ws.Range("H69").Value = 14500
$ws.Range("J69").Value = 14500
$ws.Range("L69").Value = 43500
$ws.Range("N69").Value = -45248

# Row 72: Surgical Substitution (L) | Grade 1 Mind Dissolvent
$ws.Range("H72").Value = 14500
$ws.Range("J72").Value = 14500
$ws.Range("L72").Value = 130500
$ws.Range("N72").Value = -139236

# Row 82: Rolling on Initiative | Draconian Potion of Dexterity
$ws.Range("H82").Value = 845
$ws.Range("I82").Value = 845
$ws.Range("K82").Value = 2535
$ws.Range("M82").Value = -2129

# Row 85: Darkly Dreaming Dexterity (L) | Draconian Potion of Dexterity
$ws.Range("H85").Value = 845
$ws.Range("I85").Value = 845
$ws.Range("K85").Value = 2535
$ws.Range("M85").Value = -1131

# Row 132: Fast-forwarding Flora | Growth Formula Lambda
$ws.Range("H132").Value = 2961.4
$ws.Range("I132").Value = 2961.4
$ws.Range("K132").Value = 8884.200000000001
$ws.Range("M132").Value = -6354.200000000001

# Row 137: Cutting Edge of Culinary Quality | Magnesia Whetstone
$ws.Range("H137").Value = 9797.529
$ws.Range("I137").Value = 6960.273
$ws.Range("J137").Value = 14999.167
$ws.Range("K137").Value = 20880.819
$ws.Range("L137").Value = 44997.501
$ws.Range("M137").Value = -18330.819
$ws.Range("N137").Value = -50097.501

# Row 138: All-night Crafting | Cunning Craftsman's Tisane
$ws.Range("H138").Value = 7149.7
$ws.Range("J138").Value = 7999.5713
$ws.Range("L138").Value = 23998.7139
$ws.Range("N138").Value = -34278.7139

$ws = $wb.Worksheets.Item("ARM")
# Row 122: Haste for High Durium | High Durium Nugget
$ws.Range("H122").Value = 1000
$ws.Range("I122").Value = 1000
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 3000
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -550
$ws.Range("N122").Value = ""

# Row 132: Don't Bore Me, Ore Me | Mountain Chromite Ingot
$ws.Range("H132").Value = 6825.8184
$ws.Range("I132").Value = 2760.5
$ws.Range("K132").Value = 8281.5
$ws.Range("M132").Value = -5751.5

$ws = $wb.Worksheets.Item("BSM")
# Row 139: Maul Me | Titanium Gold Maul
$ws.Range("H139").Value = 50709
$ws.Range("I139").Value = 50709
$ws.Range("K139").Value = 50709
$ws.Range("M139").Value = -45569

$ws = $wb.Worksheets.Item("CRP")
# Row 2: In with the New | Bone Harpoon
$ws.Range("H2").Value = 77
$ws.Range("J2").Value = 120
$ws.Range("L2").Value = 120
$ws.Range("N2").Value = -346

# Row 5: Bowing Out | Maple Shortbow
$ws.Range("H5").Value = 115.666664
$ws.Range("I5").Value = 79
$ws.Range("K5").Value = 79
$ws.Range("M5").Value = 33

# Row 11: Leaving without Leave | Bronze Spear
$ws.Range("H11").Value = 70
$ws.Range("I11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("M11").Value = ""

# Row 15: On the Move | Ragstone Grinding Wheel
$ws.Range("H15").Value = 760
$ws.Range("I15").Value = 700
$ws.Range("K15").Value = 700
$ws.Range("M15").Value = -530

# Row 16: Raise the Roof | Ash Lumber
$ws.Range("H16").Value = 1100
$ws.Range("I16").Value = 200
$ws.Range("K16").Value = 200
$ws.Range("M16").Value = 87

# Row 19: Shielding Sales | Square Ash Shield
$ws.Range("H19").Value = 188.5
$ws.Range("I19").Value = 188.5
$ws.Range("K19").Value = 188.5
$ws.Range("M19").Value = -18.5

# Row 24: What You Need | Square Ash Shield
$ws.Range("H24").Value = 188.5
$ws.Range("I24").Value = 188.5
$ws.Range("K24").Value = 188.5
$ws.Range("M24").Value = -18.5

# Row 31: Wall Not Found | Walnut Lumber
$ws.Range("H31").Value = 4265.5454
$ws.Range("I31").Value = 3542.3635
$ws.Range("J31").Value = 4988.727
$ws.Range("K31").Value = 3542.3635
$ws.Range("L31").Value = 4988.727
$ws.Range("M31").Value = -3247.3635
$ws.Range("N31").Value = -5578.727

# Row 34: Armoires of the Rich and Famous | Walnut Lumber
$ws.Range("H34").Value = 4265.5454
$ws.Range("I34").Value = 3542.3635
$ws.Range("J34").Value = 4988.727
$ws.Range("K34").Value = 3542.3635
$ws.Range("L34").Value = 4988.727
$ws.Range("M34").Value = -3340.3635
$ws.Range("N34").Value = -5392.727

# Row 35: Storm of Swords | Elm Macuahuitl
$ws.Range("H35").Value = 1244.3334
$ws.Range("I35").Value = 1244.3334
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 1244.3334
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -950.3334
$ws.Range("N35").Value = ""

# Row 99: O Pine | Pine Lumber
$ws.Range("H99").Value = 4649.8335
$ws.Range("I99").Value = 4649.8335
$ws.Range("K99").Value = 4649.8335
$ws.Range("M99").Value = -3151.8335

# Row 107: Built to Last | White Oak Lumber
$ws.Range("H107").Value = 3936.7896
$ws.Range("I107").Value = 4052.8823
$ws.Range("K107").Value = 4052.8823
$ws.Range("M107").Value = -2132.8823

# Row 113: Patient Patients | White Ash Lumber
$ws.Range("H113").Value = 1100
$ws.Range("I113").Value = 200
$ws.Range("K113").Value = 200
$ws.Range("M113").Value = 1970

# Row 126: A Better Conductor | Red Pine Lumber
$ws.Range("H126").Value = 4649.8335
$ws.Range("I126").Value = 4649.8335
$ws.Range("K126").Value = 13949.5005
$ws.Range("M126").Value = -11479.5005

$ws = $wb.Worksheets.Item("CUL")
# Row 2: Pork Is a Salty Food | Table Salt
$ws.Range("H2").Value = 248.57143
$ws.Range("I2").Value = 233.1
$ws.Range("J2").Value = 287.25
$ws.Range("K2").Value = 1398.6
$ws.Range("L2").Value = 1723.5
$ws.Range("M2").Value = -1285.6
$ws.Range("N2").Value = -1949.5

# Row 7: It's Always Sunny in Vylbrand | Raisins
$ws.Range("H7").Value = 63.125
$ws.Range("I7").Value = 43.57143
$ws.Range("K7").Value = 130.71429
$ws.Range("M7").Value = -18.71429000000001

# Row 11: Putting the Squeeze On | Orange Juice
$ws.Range("H11").Value = 787.5
$ws.Range("J11").Value = 716.6667
$ws.Range("L11").Value = 2150.0001
$ws.Range("N11").Value = -2430.0001

# Row 13: Fishy Revelations | Braised Pipira
$ws.Range("H13").Value = 248
$ws.Range("I13").Value = 248
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 744
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = -576
$ws.Range("N13").Value = ""

# Row 34: Fever Pitch | Chamomile Tea
$ws.Range("H34").Value = 1389.5454
$ws.Range("J34").Value = 2170.7144
$ws.Range("L34").Value = 6512.1432
$ws.Range("N34").Value = -6680.1432

# Row 81: It Goes Down Smoothly | Frozen Spirits
$ws.Range("H81").Value = 20000
$ws.Range("I81").Value = 0
$ws.Range("K81").Value = 0
$ws.Range("M81").Value = ""

# Row 84: Quenching the Flame (L) | Frozen Spirits
$ws.Range("H84").Value = 20000
$ws.Range("I84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("M84").Value = ""

# Row 121: A Cookie for Your Troubles | Coffee Biscuit
$ws.Range("H121").Value = 827.6
$ws.Range("I121").Value = 259
$ws.Range("K121").Value = 777
$ws.Range("M121").Value = 533

$ws = $wb.Worksheets.Item("GSM")
# Row 80: Needs More Prayerbell | Hardsilver Ingot
$ws.Range("H80").Value = 3208.8333
$ws.Range("I80").Value = 2922.6667
$ws.Range("J80").Value = 3495
$ws.Range("K80").Value = 2922.6667
$ws.Range("L80").Value = 3495
$ws.Range("M80").Value = -1924.6667
$ws.Range("N80").Value = -5491

# Row 83: With a Noise That Reaches Heaven (L) | Hardsilver Ingot
$ws.Range("H83").Value = 3208.8333
$ws.Range("I83").Value = 2922.6667
$ws.Range("J83").Value = 3495
$ws.Range("K83").Value = 14613.3335
$ws.Range("L83").Value = 17475
$ws.Range("M83").Value = -9621.333500000001
$ws.Range("N83").Value = -27459

# Row 102: Put the Metal to the Peddle | Durium Ingot
$ws.Range("H102").Value = 3739.4
$ws.Range("I102").Value = 3449
$ws.Range("J102").Value = 3933
$ws.Range("K102").Value = 3449
$ws.Range("L102").Value = 3933
$ws.Range("M102").Value = -1827
$ws.Range("N102").Value = -7177

# Row 113: Copious Crystal Cannons | Manasilver Nugget
$ws.Range("H113").Value = 638.6
$ws.Range("J113").Value = 582.6667
$ws.Range("L113").Value = 582.6667
$ws.Range("N113").Value = -4922.6667

# Row 122: Awarding Academic Excellence | Ametrine
$ws.Range("H122").Value = 999.3333
$ws.Range("I122").Value = 999.3333
$ws.Range("K122").Value = 2997.9999
$ws.Range("M122").Value = -547.9998999999998

# Row 123: Workplace Workout | Ametrine Ring of Fending
$ws.Range("H123").Value = 23205
$ws.Range("J123").Value = 23205
$ws.Range("L123").Value = 23205
$ws.Range("N123").Value = -28105

$ws = $wb.Worksheets.Item("LTW")
# Row 122: Hell on Leather | Gaja Leather
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").Value = ""

$ws = $wb.Worksheets.Item("WVR")
# Row 81: Where the Dragonflies, the Net Catches | Crawler Silk
$ws.Range("H81").Value = 3208
$ws.Range("J81").Value = 7500
$ws.Range("L81").Value = 15000
$ws.Range("N81").Value = -17122

# Row 84: To Kill a Dragon on Nameday (L) | Crawler Silk
$ws.Range("H84").Value = 3208
$ws.Range("J84").Value = 7500
$ws.Range("L84").Value = 75000
$ws.Range("N84").Value = -85608

# Row 107: Flax Wax | Bright Linen Yarn
$ws.Range("H107").Value = 457.72726
$ws.Range("I107").Value = 373.125
$ws.Range("K107").Value = 1119.375
$ws.Range("M107").Value = 800.625

# Row 122: Heavy Armoire | Dark Hempen Cloth
$ws.Range("H122").Value = 5000
$ws.Range("I122").Value = 5000
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 15000
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -12550
$ws.Range("N122").Value = -19900

